$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("33-3950749", "GAIL HINKEL",      94484153,  "Gail",     "Hinkel",      "05/28/1955", "87 Nelson Ave",      "Ghent",    "NY", 12075),
    @("33-3950813", "DEMETRIC FORNEY",  135643630, "Demetric", "Forney",      "09/03/1967", "367 Jerome St",      "Brooklyn", "NY", 11207),
    @("33-3950839", "ALLISON LEWIS",    132625821, "Allison",  "Lewis",       "07/08/1969", "2759 Webster Ave",   "Bronx",    "NY", 10458),
    @("33-3950858", "DOMINGO SOTO",     131506490, "Domingo",  "Soto",        "03/11/1956", "1047 Clay Ave",      "Bronx",    "NY", 10456),
    @("33-3951010", "JUDY CRUZ PALMA",  114704216, "Judy",     "Cruz Palma",  "06/18/1985", "7 Alexander Rd",     "Monroe",   "NY", 10950)
)

$startRow = 6
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]

    # Force the DOB column to stay plain text (no auto date conversion),
    # then drop back to the default "Normal" style so no stray number
    # format index is left behind on the cell.
    $dobCell = $ws.Cells.Item($r, 6)
    $dobCell.NumberFormat = "@"
    $dobCell.Value = $rowVals[5]
    $dobCell.Style = "Normal"

    $ws.Cells.Item($r, 1).Value = $rowVals[0]
    $ws.Cells.Item($r, 2).Value = $rowVals[1]
    $ws.Cells.Item($r, 3).Value = $rowVals[2]
    $ws.Cells.Item($r, 4).Value = $rowVals[3]
    $ws.Cells.Item($r, 5).Value = $rowVals[4]
    $ws.Cells.Item($r, 7).Value = $rowVals[6]
    $ws.Cells.Item($r, 8).Value = $rowVals[7]
    $ws.Cells.Item($r, 9).Value = $rowVals[8]
    $ws.Cells.Item($r, 10).Value = $rowVals[9]
}
